$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (shared strings ID -> ProjectID, Dependencies -> ProjectDependency)
$ws.Range("A1").Value = "ProjectID"
$ws.Range("C1").Value = "ProjectDependency"

# Rename the corresponding Table1 columns to match the new headers
$tbl = $ws.ListObjects.Item("Table1")
$tbl.ListColumns.Item("ID").Name = "ProjectID"
$tbl.ListColumns.Item("Dependencies").Name = "ProjectDependency"

# Move the active selection from C5 to C1
$ws.Range("C1").Select()
